$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.158.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.489.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.27%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.52%  '
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.878.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.495.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.074.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("E20").Value = '  +2.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +17.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0934'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '245.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.36%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.05'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.55%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.141'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.26%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.93'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.08'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.59%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.36'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0783'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.96'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.72%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.112'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '118.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.82%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.01%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0294'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.979.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.55%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.23%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.73%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '57.40'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.20%  '
